$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "36.507.92"
Set-TextValue "E2" "  +0.18%  "
Set-TextValue "D3" "1.954.52"
Set-TextValue "E3" "  +0.74%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "243.12"
Set-TextValue "E5" "  +0.32%  "
Set-TextValue "D6" "0.627"
Set-TextValue "E6" "  +2.84%  "
Set-TextValue "D7" "60.34"
Set-TextValue "E7" "  +6.93%  "
Set-TextValue "E8" "  -0.03%  "
Set-TextValue "E9" "  +5.24%  "
Set-TextValue "E10" "  -2.15%  "
Set-TextValue "E11" "  +0.79%  "
Set-TextValue "D12" "14.16"
Set-TextValue "E12" "  +6.66%  "
Set-TextValue "D13" "0.838"
Set-TextValue "E13" "  +4.49%  "
Set-TextValue "D14" "2.241.93"
Set-TextValue "E14" "  +0.88%  "
Set-TextValue "D15" "21.55"
Set-TextValue "E15" "  +1.21%  "
Set-TextValue "E16" "  +2.59%  "
Set-TextValue "D17" "1.956.19"
Set-TextValue "E17" "  +1.26%  "
Set-TextValue "D18" "36.481.76"
Set-TextValue "E18" "  +0.25%  "
Set-TextValue "E19" "  +0.52%  "
Set-TextValue "D20" "0.0₃0853"
Set-TextValue "E20" "  +0.10%  "
Set-TextValue "D21" "229.47"
Set-TextValue "E21" "  +1.21%  "
Set-TextValue "D22" "5.08"
Set-TextValue "E22" "  +3.11%  "
Set-TextValue "E23" "  +0.13%  "
Set-TextValue "E24" "  +2.63%  "
Set-TextValue "E25" "  +3.88%  "
Set-TextValue "D26" "0.143"
Set-TextValue "E26" "  +8.31%  "
Set-TextValue "D27" "9.17"
Set-TextValue "E27" "  +0.97%  "
Set-TextValue "D28" "160.58"
Set-TextValue "E28" "  +0.63%  "
Set-TextValue "D29" "19.29"
Set-TextValue "E29" "  +1.36%  "
Set-TextValue "D30" "1.31"
Set-TextValue "E30" "  +20.76%  "
Set-TextValue "E31" "  +2.47%  "
Set-TextValue "D32" "4.77"
Set-TextValue "E32" "  +4.79%  "
Set-TextValue "E33" "  +0.29%  "
Set-TextValue "E34" "  +7.87%  "
Set-TextValue "D35" "3.45"
Set-TextValue "E35" "  +7.82%  "
Set-TextValue "E36" "  +0.06%  "
Set-TextValue "D37" "2.26"
Set-TextValue "E37" "  +3.78%  "
Set-TextValue "D38" "1.78"
Set-TextValue "E38" "  -0.37%  "
Set-TextValue "D39" "5.45"
Set-TextValue "E39" "  -10.44%  "
Set-TextValue "D40" "0.0968"
Set-TextValue "E40" "  -2.09%  "
Set-TextValue "D41" "2.92"
Set-TextValue "E41" "  +0.43%  "
Set-TextValue "E42" "  +2.41%  "
Set-TextValue "E43" "  +1.24%  "
Set-TextValue "D44" "15.87"
Set-TextValue "E44" "  +1.26%  "
Set-TextValue "D45" "1.362.71"
Set-TextValue "E45" "  +2.62%  "
Set-TextValue "D46" "88.78"
Set-TextValue "E46" "  +3.71%  "
Set-TextValue "E47" "  +0.59%  "
Set-TextValue "E48" "  +0.97%  "
Set-TextValue "E49" "  +0.81%  "
Set-TextValue "D50" "45.21"
Set-TextValue "E50" "  +5.34%  "
Set-TextValue "D51" "2.137.34"
Set-TextValue "E51" "  +1.04%  "
